# "Page buttons fully functional"
# movies.html (sheet 3) had its data-table rows 13-34 (2nd and 3rd repeated
# "pages" worth of movie rows) cleared out, leaving just the placeholder
# column-H cells behind (blank, but still carrying their original
# formatting) wherever a page-button row used to sit, and adding two more
# blank "pages" worth of placeholder rows below for the newly working
# pagination.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)
$ws.Activate() | Out-Null

# --- existing rows 13-18 / 24-29: wipe A:G, keep H as an empty styled cell
$ws.Range("A13:G18").Clear() | Out-Null
$ws.Range("H13:H18").ClearContents() | Out-Null

$ws.Range("A24:G29").Clear() | Out-Null
$ws.Range("H24:H29").ClearContents() | Out-Null

# --- existing rows 19-23 / 30-34: these "page" rows disappear completely
$ws.Range("A19:H23").Clear() | Out-Null
$ws.Range("A30:H34").Clear() | Out-Null

# --- grab the formatting of one of the surviving placeholder cells so the
# new placeholder rows below pick up the same style (s="1") without typing
# any value into them. Re-copy before every paste batch: ClearContents()
# drops CutCopyMode, so a stale clipboard silently no-ops later pastes.
$ws.Range("H13").Copy() | Out-Null
$ws.Range("H35").PasteSpecial(-4122) | Out-Null
$ws.Range("H36").PasteSpecial(-4122) | Out-Null
$ws.Range("H37").PasteSpecial(-4122) | Out-Null
$ws.Range("H38").PasteSpecial(-4122) | Out-Null
$ws.Range("H39").PasteSpecial(-4122) | Out-Null
$ws.Range("H40").PasteSpecial(-4122) | Out-Null
$ws.Range("H35:H40").ClearContents() | Out-Null

$ws.Range("H13").Copy() | Out-Null
$ws.Range("H46").PasteSpecial(-4122) | Out-Null
$ws.Range("H47").PasteSpecial(-4122) | Out-Null
$ws.Range("H48").PasteSpecial(-4122) | Out-Null
$ws.Range("H49").PasteSpecial(-4122) | Out-Null
$ws.Range("H50").PasteSpecial(-4122) | Out-Null
$ws.Range("H51").PasteSpecial(-4122) | Out-Null
$ws.Range("H46:H51").ClearContents() | Out-Null

$excel.CutCopyMode = 0

# row 45 stays completely empty (no cells at all) but keeps a tweaked
# row height, same as it ends up with in Excel once the surrounding rows
# lose their content
$ws.Rows(45).RowHeight = 13.5

# --- view: scroll back to the top of the sheet and leave the selection on
# the first blank "page" placeholder cell
$ws.Range("F25").Select() | Out-Null

$wb.Save()
